$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits -------------------------------------------------------
# Rename the mEFCT_SPCL|... tokens to mEFCT_SPCLFUNC|... (rows 4-6, cols D/E).
# These are written before the D2 text edit so that new shared-string
# entries get appended in the same order the original author's save did.
$ws.Range("D4").Value2 = "mEFCT_SPCLFUNC|mEFCT_SHOOT"
$ws.Range("E4").Value2 = "mEFCT_SPCLFUNC|mEFCT_SHOOT"
$ws.Range("D5").Value2 = "mEFCT_SPCLFUNC|mEFCT_OPEN_BARREL"
$ws.Range("E5").Value2 = "mEFCT_SPCLFUNC|mEFCT_OPEN_BARREL"
$ws.Range("D6").Value2 = "mEFCT_SPCLFUNC|mEFCT_LOCK_LOAD"
$ws.Range("E6").Value2 = "mEFCT_SPCLFUNC|mEFCT_LOCK_LOAD"

# Shorten the introductory message in D2 (drop the trailing configure hint).
$ws.Range("D2").Value2 = "This is the FOOF Science Fiction Rubber Band Gun version 1.0."

# --- Formatting follow-through -------------------------------------------
# The shorter D2 text wraps to fewer lines, so the row shrinks.
$ws.Rows.Item(2).RowHeight = 30.75

# Column B now holds the longer "mSPCL_ONETIME | mSPCL_SHOOT" style text and
# column E was widened a bit further; set their widths to match.
$ws.Columns.Item(2).ColumnWidth = 30.66667
$ws.Columns.Item(5).ColumnWidth = 38.66667

# --- View state ------------------------------------------------------------
# Scroll back up to the top of the sheet and select C4.
$ws.Range("C4").Select()
